$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- X7 / Y7 (Saturday), X8 / Y8 (Sunday) ---------------------------------
# Copy number/format style from the existing "week" rows (X3/Y3) so the new
# cells pick up style indexes 14 (centered date) and 9 (centered text).
$ws.Range("X3").Copy()
$ws.Range("X7").PasteSpecial(-4122)
$ws.Range("X7").Value = 46011

$ws.Range("Y3").Copy()
$ws.Range("Y7").PasteSpecial(-4122)
$ws.Range("Y7").Value = "Saturday"

$ws.Range("X3").Copy()
$ws.Range("X8").PasteSpecial(-4122)
$ws.Range("X8").Value = 46012

$ws.Range("Y3").Copy()
$ws.Range("Y8").PasteSpecial(-4122)
$ws.Range("Y8").Value = "Sunday"

# --- E23: new Github link hyperlink + text --------------------------------
$e23 = $ws.Range("E23")
$ws.Hyperlinks.Add($e23, "https://github.com/AakashChidambaranathan/Intership_task/tree/8702c10b09067fcb834f14066909a04ddc83680e", [System.Type]::Missing, [System.Type]::Missing, "AakashChidambaranathan/Intership_task at 8702c10b09067fcb834f14066909a04ddc83680e")
# Restore the normal hyperlink-cell look (style used by the other E-column
# hyperlink cells) since Hyperlinks.Add generates its own fresh style.
$ws.Range("E19").Copy()
$e23.PasteSpecial(-4122)

# --- F23: task system location path ---------------------------------------
$ws.Range("F23").Value = "D:\intership\task\third_week\task_3_Blog_on_react\Intership_task\third_week\task_3_Blog_on_react\blog-app"

# --- H19: weekly summary text ----------------------------------------------
$ws.Range("H3").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H19").Value = "During this week, I developed and enhanced a blog application using React and Bootstrap by adding popups, animations, responsive design, interactive features, and dynamic text updates. I also learned and worked with Node.js and Express, including middleware concepts and storing application data in JSON files, with hands-on practice and work-from-home tasks."

# --- H20:H23 restyle (still part of the H19:H23 merge) ----------------------
$ws.Range("H12").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H23").PasteSpecial(-4122)

# --- sheet view: zoom + selection -------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 110
$ws.Range("C12").Select()
